# Auto-generated edit script: applies cell-level value updates to the FFXIV leve-profit workbook
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, matching the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 7158957
$ws.Range("I32").Value = 6999.4
$ws.Range("J32").Value = 11132266
$ws.Range("K32").Value = 6999.4
$ws.Range("L32").Value = 11132266
$ws.Range("M32").Value = -6673.4
$ws.Range("N32").Value = -11132918
$ws.Range("H55").Value = 190.72728
$ws.Range("I55").Value = 99.8
$ws.Range("J55").Value = 266.5
$ws.Range("K55").Value = 99.8
$ws.Range("L55").Value = 266.5
$ws.Range("M55").Value = 114.2
$ws.Range("N55").Value = -694.5
$ws.Range("H70").Value = 1628.6923
$ws.Range("I70").Value = 1874
$ws.Range("J70").Value = 1555.1
$ws.Range("K70").Value = 5622
$ws.Range("L70").Value = 4665.299999999999
$ws.Range("M70").Value = -5352
$ws.Range("N70").Value = -5205.299999999999
$ws.Range("H73").Value = 1628.6923
$ws.Range("I73").Value = 1874
$ws.Range("J73").Value = 1555.1
$ws.Range("K73").Value = 5622
$ws.Range("L73").Value = 4665.299999999999
$ws.Range("M73").Value = -4686
$ws.Range("N73").Value = -6537.299999999999
$ws.Range("H113").Value = 4379.7
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 4588.5557
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 4588.5557
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -11096.5557
$ws.Range("H125").Value = 5031.091
$ws.Range("J125").Value = 6663.1665
$ws.Range("L125").Value = 59968.4985
$ws.Range("N125").Value = -64888.4985
$ws.Range("H134").Value = 99829.164
$ws.Range("J134").Value = 99829.164
$ws.Range("L134").Value = 99829.164
$ws.Range("N134").Value = -109969.164

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 453
$ws.Range("I16").Value = 453
$ws.Range("K16").Value = 453
$ws.Range("M16").Value = -166
$ws.Range("H19").Value = 608
$ws.Range("I19").Value = 608
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 608
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -379
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 6509.769
$ws.Range("I32").Value = 2250.535
$ws.Range("K32").Value = 2250.535
$ws.Range("M32").Value = -1963.535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5379.619
$ws.Range("I134").Value = 3032.7856
$ws.Range("J134").Value = 10073.286
$ws.Range("K134").Value = 9098.356800000001
$ws.Range("L134").Value = 30219.858
$ws.Range("M134").Value = -6563.356800000001
$ws.Range("N134").Value = -35289.858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 381.14285
$ws.Range("I22").Value = 233.6
$ws.Range("J22").Value = 750
$ws.Range("K22").Value = 233.6
$ws.Range("L22").Value = 750
$ws.Range("M22").Value = 116.4
$ws.Range("N22").Value = -1450
$ws.Range("H25").Value = 538
$ws.Range("I25").Value = 300.5
$ws.Range("K25").Value = 300.5
$ws.Range("M25").Value = -126.5
$ws.Range("H62").Value = 4399.8
$ws.Range("I62").Value = 4499.75
$ws.Range("K62").Value = 4499.75
$ws.Range("M62").Value = -3875.75
$ws.Range("H65").Value = 4399.8
$ws.Range("I65").Value = 4499.75
$ws.Range("K65").Value = 22498.75
$ws.Range("M65").Value = -19378.75
$ws.Range("H99").Value = 6781679
$ws.Range("I99").Value = 7409964.5
$ws.Range("J99").Value = 5210966.5
$ws.Range("K99").Value = 7409964.5
$ws.Range("L99").Value = 5210966.5
$ws.Range("M99").Value = -7408466.5
$ws.Range("N99").Value = -5213962.5
$ws.Range("H107").Value = 1809.9584
$ws.Range("I107").Value = 1215.375
$ws.Range("J107").Value = 2999.125
$ws.Range("K107").Value = 1215.375
$ws.Range("L107").Value = 2999.125
$ws.Range("M107").Value = 704.625
$ws.Range("N107").Value = -6839.125
$ws.Range("H126").Value = 6781679
$ws.Range("I126").Value = 7409964.5
$ws.Range("J126").Value = 5210966.5
$ws.Range("K126").Value = 22229893.5
$ws.Range("L126").Value = 15632899.5
$ws.Range("M126").Value = -22227423.5
$ws.Range("N126").Value = -15637839.5
$ws.Range("H132").Value = 407869.25
$ws.Range("I132").Value = 1573.3214
$ws.Range("J132").Value = 3251940.8
$ws.Range("K132").Value = 4719.9642
$ws.Range("L132").Value = 9755822.399999999
$ws.Range("M132").Value = -2189.9642
$ws.Range("N132").Value = -9760882.399999999
$ws.Range("H134").Value = 75198.36
$ws.Range("I134").Value = 4603
$ws.Range("J134").Value = 169325.5
$ws.Range("K134").Value = 13809
$ws.Range("L134").Value = 507976.5
$ws.Range("M134").Value = -11274
$ws.Range("N134").Value = -513046.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 116.28571
$ws.Range("I2").Value = 258.6
$ws.Range("J2").Value = 37.22222
$ws.Range("K2").Value = 1551.6
$ws.Range("L2").Value = 223.33332
$ws.Range("M2").Value = -1438.6
$ws.Range("N2").Value = -449.33332
$ws.Range("H5").Value = 599.2
$ws.Range("I5").Value = 385.0909
$ws.Range("J5").Value = 723.1579
$ws.Range("K5").Value = 1155.2727
$ws.Range("L5").Value = 2169.4737
$ws.Range("M5").Value = -1043.2727
$ws.Range("N5").Value = -2393.4737
$ws.Range("H26").Value = 100929
$ws.Range("J26").Value = 143777.14
$ws.Range("L26").Value = 431331.42
$ws.Range("N26").Value = -431907.42
$ws.Range("H56").Value = 6248.75
$ws.Range("I56").Value = 6248.75
$ws.Range("K56").Value = 6248.75
$ws.Range("M56").Value = -5718.75
$ws.Range("H59").Value = 3573.75
$ws.Range("I59").Value = 3033
$ws.Range("K59").Value = 9099
$ws.Range("M59").Value = -8559
$ws.Range("H60").Value = 589.25
$ws.Range("I60").Value = 34.166668
$ws.Range("J60").Value = 922.3
$ws.Range("K60").Value = 102.500004
$ws.Range("L60").Value = 2766.9
$ws.Range("M60").Value = 148.499996
$ws.Range("N60").Value = -3268.9
$ws.Range("H93").Value = 400
$ws.Range("I93").Value = 400
$ws.Range("K93").Value = 1200
$ws.Range("M93").Value = 672
$ws.Range("H98").Value = 10499.333
$ws.Range("J98").Value = 12498.5
$ws.Range("L98").Value = 37495.5
$ws.Range("N98").Value = -40491.5
$ws.Range("H104").Value = 4793.05
$ws.Range("J104").Value = 4992.1055
$ws.Range("L104").Value = 14976.3165
$ws.Range("N104").Value = -20218.3165
$ws.Range("H113").Value = 42502.918
$ws.Range("I113").Value = 1149.8
$ws.Range("K113").Value = 3449.4
$ws.Range("M113").Value = -1279.4
$ws.Range("H131").Value = 1691.125
$ws.Range("I131").Value = 1359.8
$ws.Range("K131").Value = 4079.4
$ws.Range("M131").Value = 960.6000000000004
$ws.Range("H132").Value = 7181.933
$ws.Range("J132").Value = 7181.933
$ws.Range("L132").Value = 64637.397
$ws.Range("N132").Value = -69697.397
$ws.Range("H135").Value = 599.2
$ws.Range("I135").Value = 385.0909
$ws.Range("J135").Value = 723.1579
$ws.Range("K135").Value = 3465.8181
$ws.Range("L135").Value = 6508.4211
$ws.Range("M135").Value = -930.8181
$ws.Range("N135").Value = -11578.4211

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 2502.5
$ws.Range("I18").Value = 2003.3334
$ws.Range("K18").Value = 2003.3334
$ws.Range("M18").Value = -1710.3334
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H102").Value = 1265.8572
$ws.Range("I102").Value = 1277.3334
$ws.Range("K102").Value = 1277.3334
$ws.Range("M102").Value = 344.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 74500
$ws.Range("J6").Value = 74500
$ws.Range("L6").Value = 74500
$ws.Range("N6").Value = -74724
$ws.Range("H22").Value = 8979.385
$ws.Range("I22").Value = 1448.3636
$ws.Range("J22").Value = 50400
$ws.Range("K22").Value = 1448.3636
$ws.Range("L22").Value = 50400
$ws.Range("M22").Value = -1153.3636
$ws.Range("N22").Value = -50990
$ws.Range("H23").Value = 2000
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = -1770
$ws.Range("N23").Value = -2460
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()
$ws.Range("H27").Value = 8979.385
$ws.Range("I27").Value = 1448.3636
$ws.Range("J27").Value = 50400
$ws.Range("K27").Value = 1448.3636
$ws.Range("L27").Value = 50400
$ws.Range("M27").Value = -1341.3636
$ws.Range("N27").Value = -50614
$ws.Range("H120").Value = 35000
$ws.Range("J120").Value = 35000
$ws.Range("L120").Value = 35000
$ws.Range("N120").Value = -44676
$ws.Range("H122").Value = 92860960
$ws.Range("J122").Value = 42861060
$ws.Range("L122").Value = 128583180
$ws.Range("N122").Value = -128588080
$ws.Range("H128").Value = 69996.5
$ws.Range("J128").Value = 69996.5
$ws.Range("L128").Value = 69996.5
$ws.Range("N128").Value = -79956.5
$ws.Range("H136").Value = 3161.7144
$ws.Range("I136").Value = 3305.7827
$ws.Range("K136").Value = 9917.348100000001
$ws.Range("M136").Value = -7367.348100000001
$ws.Range("H137").Value = 109995
$ws.Range("J137").Value = 109995
$ws.Range("L137").Value = 109995
$ws.Range("N137").Value = -120195

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H116").Value = 77425
$ws.Range("J116").Value = 77425
$ws.Range("L116").Value = 77425
$ws.Range("N116").Value = -86603
$ws.Range("H132").Value = 2024.6487
$ws.Range("I132").Value = 1776.3478
$ws.Range("K132").Value = 5379.0434
$ws.Range("M132").Value = -2799.0434
